$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$tcs2 = $nm.Theme.ThemeColorScheme
$tcs2.Item(3).RGB = 999
Write-Host "done"
